$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.285.90'
$ws.Range('E2').Value = '  -1.04%  '
$ws.Range('D3').Value = '2.479.68'
$ws.Range('E3').Value = '  -0.74%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '519.66'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.95'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.558'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.56%  '
$ws.Range('D9').Value = '2.493.54'
$ws.Range('E9').Value = '  -0.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0988'
$ws.Range('D10').ClearFormats()
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.156'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.79%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.32'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.07%  '
$ws.Range('E13').Value = '  -1.96%  '
$ws.Range('D14').Value = '2.920.93'
$ws.Range('E14').Value = '  -0.63%  '
$ws.Range('D15').Value = '58.187.74'
$ws.Range('E15').Value = '  -1.07%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.17'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -2.20%  '
$ws.Range('E17').Value = '  -2.00%  '
$ws.Range('D18').Value = '2.494.95'
$ws.Range('E18').Value = '  -0.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.67'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -3.19%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.18'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '320.75'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.52%  '
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.75'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -3.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.42'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.89%  '
$ws.Range('E25').Value = '  -1.84%  '
$ws.Range('B26').Value = 'Binance-PegBSC-USD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.60%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.162'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.38'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.60%  '
$ws.Range('D29').Value = '0.0₃0750'
$ws.Range('E29').Value = '  -1.10%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '169.63'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.31%  '
$ws.Range('B31').Value = 'Aptos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.32'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.89%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.69'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.73%  '
$ws.Range('E33').Value = '  +2.23%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('E35').Value = '  -0.24%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.13'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.04%  '
$ws.Range('E37').Value = '  -0.75%  '
$ws.Range('E38').Value = '  +0.12%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.68'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.22%  '
$ws.Range('E40').Value = '  -3.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.797'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.14'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +3.27%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '275.30'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.05%  '
$ws.Range('E44').Value = '  -2.93%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.598'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '123.99'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -4.26%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0909'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.57%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0491'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.38%  '
$ws.Range('E49').Value = '  -1.44%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.06'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.81%  '
$ws.Range('D51').Value = '1.740.06'
$ws.Range('E51').Value = '  -0.70%  '
